# Applies the value updates described in the commit diff to the
# "Bahamut_Profits" workbook (sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

# ALC!61
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 300
$ws.Range("I61").Value = 300
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 900
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -728
$ws.Range("N61").ClearContents()

# ALC!76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 23917230
$ws.Range("I76").Value = 39289600
$ws.Range("J76").Value = 4654.278
$ws.Range("K76").Value = 39289600
$ws.Range("L76").Value = 4654.278
$ws.Range("M76").Value = -39289285
$ws.Range("N76").Value = -5284.278

# ALC!79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 23917230
$ws.Range("I79").Value = 39289600
$ws.Range("J79").Value = 4654.278
$ws.Range("K79").Value = 39289600
$ws.Range("L79").Value = 4654.278
$ws.Range("M79").Value = -39288508
$ws.Range("N79").Value = -6838.278

# ALC!106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3683.1667
$ws.Range("I106").Value = 3724.75
$ws.Range("J106").Value = 3600
$ws.Range("K106").Value = 3724.75
$ws.Range("L106").Value = 3600
$ws.Range("M106").Value = -3093.75
$ws.Range("N106").Value = -4862

# ALC!107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1666933.5
$ws.Range("I107").Value = 2500112.8
$ws.Range("J107").Value = 575
$ws.Range("K107").Value = 2500112.8
$ws.Range("L107").Value = 575
$ws.Range("M107").Value = -2498192.8
$ws.Range("N107").Value = -4415

# ARM!64
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 29000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 29000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 29000
$ws.Range("N64").Value = -29496

# ARM!67
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 29000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 29000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 29000
$ws.Range("N67").Value = -30716

# ARM!110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 836.6667
$ws.Range("I110").Value = 836.6667
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 836.6667
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1208.3333

# BSM!88
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 40000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 40000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 40000
$ws.Range("N88").Value = -40812

# BSM!91
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 40000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 40000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 40000
$ws.Range("N91").Value = -42808

# BSM!105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2965.9412
$ws.Range("I105").Value = 1978.8889
$ws.Range("J105").Value = 4076.375
$ws.Range("K105").Value = 1978.8889
$ws.Range("L105").Value = 4076.375
$ws.Range("M105").Value = -231.8888999999999
$ws.Range("N105").Value = -7570.375

# BSM!107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7180.7617
$ws.Range("I107").Value = 882.17645
$ws.Range("J107").Value = 33949.75
$ws.Range("K107").Value = 882.17645
$ws.Range("L107").Value = 33949.75
$ws.Range("M107").Value = 1037.82355
$ws.Range("N107").Value = -37789.75

# CRP!62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 9385.714
$ws.Range("I62").Value = 9933.333
$ws.Range("J62").Value = 8400
$ws.Range("K62").Value = 9933.333
$ws.Range("L62").Value = 8400
$ws.Range("M62").Value = -9309.333
$ws.Range("N62").Value = -9648

# CRP!65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 9385.714
$ws.Range("I65").Value = 9933.333
$ws.Range("J65").Value = 8400
$ws.Range("K65").Value = 49666.665
$ws.Range("L65").Value = 42000
$ws.Range("M65").Value = -46546.665
$ws.Range("N65").Value = -48240

# CRP!99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2750
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -7996

# CRP!107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 324.07895
$ws.Range("I107").Value = 181.8
$ws.Range("J107").Value = 374.89285
$ws.Range("K107").Value = 181.8
$ws.Range("L107").Value = 374.89285
$ws.Range("M107").Value = 1738.2
$ws.Range("N107").Value = -4214.89285

# CRP!126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2750
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -19940

# CUL!54
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 600
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 600
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 1800
$ws.Range("N54").Value = -2918

# CUL!60
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 191.83333
$ws.Range("I60").Value = 114.75
$ws.Range("J60").Value = 346
$ws.Range("K60").Value = 344.25
$ws.Range("L60").Value = 1038
$ws.Range("M60").Value = -93.25
$ws.Range("N60").Value = -1540

# CUL!61
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 195
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 195
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 585
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -1015

# GSM!70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4662.8823
$ws.Range("I70").Value = 4100
$ws.Range("J70").Value = 4783.5
$ws.Range("K70").Value = 4100
$ws.Range("L70").Value = 4783.5
$ws.Range("M70").Value = -3830
$ws.Range("N70").Value = -5323.5

# GSM!73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4662.8823
$ws.Range("I73").Value = 4100
$ws.Range("J73").Value = 4783.5
$ws.Range("K73").Value = 4100
$ws.Range("L73").Value = 4783.5
$ws.Range("M73").Value = -3164
$ws.Range("N73").Value = -6655.5

# GSM!80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3420
$ws.Range("I80").Value = 3409
$ws.Range("J80").Value = 3475
$ws.Range("K80").Value = 3409
$ws.Range("L80").Value = 3475
$ws.Range("M80").Value = -2411
$ws.Range("N80").Value = -5471

# GSM!83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3420
$ws.Range("I83").Value = 3409
$ws.Range("J83").Value = 3475
$ws.Range("K83").Value = 17045
$ws.Range("L83").Value = 17375
$ws.Range("M83").Value = -12053
$ws.Range("N83").Value = -27359

# GSM!87
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").ClearContents()

# GSM!90
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").ClearContents()

# GSM!102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1290.5
$ws.Range("I102").Value = 887.3333
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 887.3333
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = 734.6667
$ws.Range("N102").Value = -5744

# LTW!70
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# LTW!73
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# LTW!82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2186.6667
$ws.Range("I82").Value = 2314.2856
$ws.Range("J82").Value = 2075
$ws.Range("K82").Value = 2314.2856
$ws.Range("L82").Value = 2075
$ws.Range("M82").Value = -1953.2856
$ws.Range("N82").Value = -2797

# LTW!85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2186.6667
$ws.Range("I85").Value = 2314.2856
$ws.Range("J85").Value = 2075
$ws.Range("K85").Value = 2314.2856
$ws.Range("L85").Value = 2075
$ws.Range("M85").Value = -1066.2856
$ws.Range("N85").Value = -4571

# LTW!88
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

# LTW!91
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

# LTW!122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5184.4194
$ws.Range("I122").Value = 5549.148
$ws.Range("J122").Value = 2722.5
$ws.Range("K122").Value = 16647.444
$ws.Range("L122").Value = 8167.5
$ws.Range("M122").Value = -14197.444
$ws.Range("N122").Value = -13067.5
